{"js": "// Add a new closing paragraph at the very end of the document body\n// (right before the final section break), reading:\n//   \"Finished the Lesson 3, manipulating variables, expressions and statements!\"\n//\n// The source OOXML splits this sentence across three runs (identical\n// formatting on each: 12pt / en-US), so we build that exact paragraph as a\n// small Flat-OPC OOXML fragment and insert it with insertOoxml(). Using the\n// higher-level insertParagraph()/insertText() calls would leave the correct\n// text but would let the host coalesce the adjacent, identically-formatted\n// runs into a single <w:r> - inserting ready-made OOXML keeps the three run\n// boundaries intact.\n\nconst paragraphOoxml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\">Finished the Lesson 3, </w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t>manipulating</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"24\"/>\n                <w:szCs w:val=\"24\"/>\n                <w:lang w:val=\"en-US\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> variables, expressions and statements!</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst body = context.document.body;\n\n// Create a fresh empty paragraph at the end of the body, then replace its\n// (empty) contents with our prebuilt OOXML paragraph so the run structure\n// above is preserved exactly instead of being normalised/merged.\nconst newPara = body.insertParagraph(\"\", Word.InsertLocation.end);\nnewPara.insertOoxml(paragraphOoxml, Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Add a new closing paragraph at the very end of the document body (right\n# before the final section break), reading:\n#   \"Finished the Lesson 3, manipulating variables, expressions and statements!\"\n#\n# The source OOXML splits this sentence across three runs that all share the\n# same formatting (12pt / en-US). Building the paragraph up with\n# Range.InsertAfter()/.Text would leave the wording correct but would let the\n# host coalesce adjacent, identically-formatted runs into a single <w:r>, so\n# instead we hand the exact paragraph markup to Range.InsertXML(), which\n# inserts it verbatim (three separate <w:r> elements) at the collapsed end of\n# the document.\n\n$d = $word.ActiveDocument\n\n$end = $d.Content\n$end.Collapse(0)\n\n$paragraphXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p>' +\n              '<w:pPr><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr>' +\n              '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\">Finished the Lesson 3, </w:t></w:r>' +\n              '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>manipulating</w:t></w:r>' +\n              '<w:r><w:rPr><w:sz w:val=\"24\"/><w:szCs w:val=\"24\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> variables, expressions and statements!</w:t></w:r>' +\n            '</w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n$end.InsertXML($paragraphXml)\n"}
